# check #11 on 31/3/2025
# Remove specific rows from the "tools difference" report (Sheet1).
# Rows to delete (original row numbers, before any deletion):
#   28 -> Kitagwenda / Muyenga / Muyenga_Central / Spray Pumps
#   61 -> Rakai / Kanoni / Kanoni_Rakai / Mortar and Pestle
#   80 -> Rakai / Kasensero / Nabyala / Filters
#   81 -> Rakai / Kasensero / Nabyala / Mortar and Pestle
#   84 -> Rakai / Kasensero / Nabyala / Watering can
#
# Delete from the bottom up so earlier row numbers stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowsToDelete = @(84, 81, 80, 61, 28)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
